# Plantilla de Casos de Uso — CU 18 / CU 19 status update
# CU-18 "Generar reporte de ingresos y egresos" (row 22) and
# CU-19 "CRU renta espacio" (row 23) move from Estado "vacio" to
# "planificado" and get 1 hr of Esfuerzo logged, now that they've been
# picked up (robustness/sequence diagrams, descriptions, domain model).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# CU- 18 : Generar reporte de ingresos y egresos
$ws.Range("E22").Value = "planificado"
$ws.Range("F22").Value = 1

# CU- 19 : CRU renta espacio
$ws.Range("E23").Value = "planificado"
$ws.Range("F23").Value = 1

# Match the author's saved viewport/selection state
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("F23").Select()
